# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" quarter sheet (positioned right after "总计"
# and before the existing "2022-Q3" sheet) and records the new quarter's
# totals on the "总计" summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary ("总计") sheet: shift the existing two rows down and add
#    the new 2022-Q4 row on top, newest quarter first.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Row 4 <- old row 3 (2021-Q3)
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(4, 2).Value = "2021-Q3"
$summary.Cells.Item(4, 3).Value = 1
$summary.Cells.Item(4, 4).Value = 0.08

# Row 3 <- old row 2 (2022-Q3)
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(3, 2).Value = "2022-Q3"
$summary.Cells.Item(3, 3).Value = 5
$summary.Cells.Item(3, 4).Value = 0.08

# Row 2 <- new 2022-Q4 figures
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 6
$summary.Cells.Item(2, 4).Value = 0.15

# Row 4's index cell (A4) is brand new territory on the sheet; match the
# look of the existing index column (A2/A3) by copying A3's formatting.
$summary.Cells.Item(3, 1).Copy($summary.Cells.Item(4, 1))
$summary.Cells.Item(4, 1).Value = 2

# ---------------------------------------------------------------------
# 2. Duplicate the existing "2022-Q3" sheet (keeps header/number
#    formatting intact), rename it to "2022-Q4" and place it right
#    before "2022-Q3" so the tab order reads 总计, 2022-Q4, 2022-Q3,
#    2021-Q3.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The source has 5 fund rows (2-6); 2022-Q4 needs a 6th (row 7) -- copy
# row 6's formatting down so the new row matches the sheet's style.
$q4.Range("A6:H6").Copy($q4.Range("A7"))

# ---------------------------------------------------------------------
# 3. Fill in the 2022-Q4 fund holdings.
#    Column A is the 0-based row index and column H is the numeric
#    rank; B (fund code) and D:G (size/position percentages) are
#    numeric-looking text in the source data (fund codes keep leading
#    zeros), so mark them as Text before writing so Excel doesn't
#    silently convert them to numbers. C (fund name) is never
#    numeric-looking, so it needs no special handling.
# ---------------------------------------------------------------------
$q4.Range("B2:B7").NumberFormat = "@"
$q4.Range("D2:G7").NumberFormat = "@"

$rows = @(
    @(0, "952035", "国泰君安君得诚混合",               "2.21", "85.83", "3.04", "0.0672", 10),
    @(1, "161224", "国投瑞银新丝路灵活配置混合（LOF）", "0.84", "93.28", "4.37", "0.0367", 9),
    @(2, "012432", "国投瑞银安泰混合C",                "1.27", "35.74", "2.16", "0.0274", 7),
    @(3, "012019", "国投瑞银安泽混合A",                "0.62", "31.81", "2.11", "0.0131", 5),
    @(4, "012020", "国投瑞银安泽混合C",                "0.11", "31.81", "2.11", "0.0023", 5),
    @(5, "012431", "国投瑞银安泰混合A",                "0.10", "35.74", "2.16", "0.0022", 7)
)

foreach ($row in $rows) {
    $r = [int]$row[0] + 2
    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = $row[3]
    $q4.Cells.Item($r, 5).Value = $row[4]
    $q4.Cells.Item($r, 6).Value = $row[5]
    $q4.Cells.Item($r, 7).Value = $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
}
